# Update "想去人数" (column F) figures across all four sheets to match the
# regenerated gh-pages data snapshot (commit 456a3b4).

$wb = $excel.ActiveWorkbook

# 展览
$ws = $wb.Worksheets.Item("展览")
$ws.Cells.Item(3, 6).Value = 344
$ws.Cells.Item(4, 6).Value = 395
$ws.Cells.Item(5, 6).Value = 88
$ws.Cells.Item(8, 6).Value = 755
$ws.Cells.Item(10, 6).Value = 5946
$ws.Cells.Item(12, 6).Value = 1691
$ws.Cells.Item(14, 6).Value = 5701
$ws.Cells.Item(15, 6).Value = 5701
$ws.Cells.Item(16, 6).Value = 103
$ws.Cells.Item(18, 6).Value = 146
$ws.Cells.Item(20, 6).Value = 1603
$ws.Cells.Item(23, 6).Value = 130
$ws.Cells.Item(24, 6).Value = 1259
$ws.Cells.Item(27, 6).Value = 21

# 演出
$ws = $wb.Worksheets.Item("演出")
$ws.Cells.Item(4, 6).Value = 124
$ws.Cells.Item(5, 6).Value = 253
$ws.Cells.Item(8, 6).Value = 334
$ws.Cells.Item(15, 6).Value = 86
$ws.Cells.Item(22, 6).Value = 1

# 本地生活
$ws = $wb.Worksheets.Item("本地生活")
$ws.Cells.Item(2, 6).Value = 9473
$ws.Cells.Item(3, 6).Value = 2198
$ws.Cells.Item(4, 6).Value = 573

# 全部类型
$ws = $wb.Worksheets.Item("全部类型")
$ws.Cells.Item(2, 6).Value = 9473
$ws.Cells.Item(3, 6).Value = 2198
$ws.Cells.Item(4, 6).Value = 573
$ws.Cells.Item(5, 6).Value = 344
$ws.Cells.Item(6, 6).Value = 395
$ws.Cells.Item(7, 6).Value = 88
$ws.Cells.Item(12, 6).Value = 124
$ws.Cells.Item(13, 6).Value = 253
$ws.Cells.Item(14, 6).Value = 755
$ws.Cells.Item(16, 6).Value = 5946
$ws.Cells.Item(18, 6).Value = 335
$ws.Cells.Item(19, 6).Value = 1691
$ws.Cells.Item(25, 6).Value = 5701
$ws.Cells.Item(26, 6).Value = 5701
$ws.Cells.Item(27, 6).Value = 103
$ws.Cells.Item(29, 6).Value = 146
$ws.Cells.Item(31, 6).Value = 1603
$ws.Cells.Item(34, 6).Value = 130
$ws.Cells.Item(35, 6).Value = 1259
$ws.Cells.Item(40, 6).Value = 86
$ws.Cells.Item(41, 6).Value = 21
$ws.Cells.Item(50, 6).Value = 1
